# Auto-generated edit script applying the Faerie Profits value updates
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 326467
$ws.Range("J17").Value = 326467
$ws.Range("L17").Value = 979401
$ws.Range("N17").Value = -979737
$ws.Range("H99").Value = 1504.3334
$ws.Range("I99").Value = 2158.6
$ws.Range("J99").Value = 686.5
$ws.Range("K99").Value = 6475.799999999999
$ws.Range("L99").Value = 2059.5
$ws.Range("M99").Value = -4977.799999999999
$ws.Range("N99").Value = -5055.5
$ws.Range("H111").Value = 12563.346
$ws.Range("I111").Value = 13990.333
$ws.Range("J111").Value = 6570
$ws.Range("K111").Value = 41970.999
$ws.Range("L111").Value = 19710
$ws.Range("M111").Value = -38903.999
$ws.Range("N111").Value = -25844
$ws.Range("H129").Value = 76928000
$ws.Range("I129").Value = 125000800
$ws.Range("J129").Value = 11523.4
$ws.Range("K129").Value = 375002400
$ws.Range("L129").Value = 34570.2
$ws.Range("M129").Value = -374997400
$ws.Range("N129").Value = -44570.2
$ws.Range("H137").Value = 2231.8215
$ws.Range("I137").Value = 2150.7073
$ws.Range("K137").Value = 6452.1219
$ws.Range("M137").Value = -3902.1219
$ws.Range("H138").Value = 104444.42
$ws.Range("I138").Value = 1191.091
$ws.Range("J138").Value = 117351.09
$ws.Range("K138").Value = 3573.273
$ws.Range("L138").Value = 352053.27
$ws.Range("M138").Value = 1566.727
$ws.Range("N138").Value = -362333.27
$ws.Range("H140").Value = 94814.164
$ws.Range("J140").Value = 94814.164
$ws.Range("L140").Value = 94814.164
$ws.Range("N140").Value = -105174.164

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value = 12999
$ws.Range("I36").Value = 12999
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 12999
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -12653
$ws.Range("N36").ClearContents()
$ws.Range("H45").Value = 2798.1035
$ws.Range("I45").Value = 3224.4614
$ws.Range("K45").Value = 3224.4614
$ws.Range("M45").Value = -2847.4614
$ws.Range("H49").Value = 19999
$ws.Range("J49").Value = 19999
$ws.Range("L49").Value = 19999
$ws.Range("N49").Value = -20519
$ws.Range("H61").Value = 11250.219
$ws.Range("I61").Value = 6340.44
$ws.Range("J61").Value = 28785.143
$ws.Range("K61").Value = 6340.44
$ws.Range("L61").Value = 28785.143
$ws.Range("M61").Value = -6128.44
$ws.Range("N61").Value = -29209.143
$ws.Range("H74").Value = 4667.4375
$ws.Range("I74").Value = 1947.75
$ws.Range("J74").Value = 7387.125
$ws.Range("K74").Value = 1947.75
$ws.Range("L74").Value = 7387.125
$ws.Range("M74").Value = -1073.75
$ws.Range("N74").Value = -9135.125
$ws.Range("H77").Value = 4667.4375
$ws.Range("I77").Value = 1947.75
$ws.Range("J77").Value = 7387.125
$ws.Range("K77").Value = 9738.75
$ws.Range("L77").Value = 36935.625
$ws.Range("M77").Value = -5370.75
$ws.Range("N77").Value = -45671.625
$ws.Range("H122").Value = 2746.2104
$ws.Range("I122").Value = 2292.4194
$ws.Range("K122").Value = 6877.2582
$ws.Range("M122").Value = -4427.2582
$ws.Range("H136").Value = 11250.219
$ws.Range("I136").Value = 6340.44
$ws.Range("J136").Value = 28785.143
$ws.Range("K136").Value = 19021.32
$ws.Range("L136").Value = 86355.429
$ws.Range("M136").Value = -16471.32
$ws.Range("N136").Value = -91455.429

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2190.3635
$ws.Range("I20").Value = 2214.5
$ws.Range("K20").Value = 2214.5
$ws.Range("M20").Value = -1967.5
$ws.Range("H22").Value = 307.55173
$ws.Range("I22").Value = 325.7619
$ws.Range("K22").Value = 325.7619
$ws.Range("M22").Value = -152.7619
$ws.Range("H54").Value = 3998.5
$ws.Range("I54").Value = 1687.2
$ws.Range("J54").Value = 15555
$ws.Range("K54").Value = 1687.2
$ws.Range("L54").Value = 15555
$ws.Range("M54").Value = -1203.2
$ws.Range("N54").Value = -16523
$ws.Range("H105").Value = 6539.364
$ws.Range("I105").Value = 5858
$ws.Range("K105").Value = 5858
$ws.Range("M105").Value = -4111
$ws.Range("H107").Value = 847.8077
$ws.Range("I107").Value = 787.3043
$ws.Range("K107").Value = 787.3043
$ws.Range("M107").Value = 1132.6957

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2236.9583
$ws.Range("I31").Value = 1573.8182
$ws.Range("K31").Value = 1573.8182
$ws.Range("M31").Value = -1278.8182
$ws.Range("H34").Value = 2236.9583
$ws.Range("I34").Value = 1573.8182
$ws.Range("K34").Value = 1573.8182
$ws.Range("M34").Value = -1371.8182
$ws.Range("H58").Value = 2906.4092
$ws.Range("I58").Value = 2763.8462
$ws.Range("J58").Value = 3112.3333
$ws.Range("K58").Value = 2763.8462
$ws.Range("L58").Value = 3112.3333
$ws.Range("M58").Value = -2560.8462
$ws.Range("N58").Value = -3518.3333
$ws.Range("H122").Value = 3892.5356
$ws.Range("I122").Value = 3748.682
$ws.Range("K122").Value = 11246.046
$ws.Range("M122").Value = -8796.045999999998
$ws.Range("H132").Value = 1741836.9
$ws.Range("I132").Value = 3079246.2
$ws.Range("J132").Value = 3204.7
$ws.Range("K132").Value = 9237738.600000001
$ws.Range("L132").Value = 9614.099999999999
$ws.Range("M132").Value = -9235208.600000001
$ws.Range("N132").Value = -14674.1
$ws.Range("H134").Value = 3705.75
$ws.Range("I134").Value = 2257.25
$ws.Range("K134").Value = 6771.75
$ws.Range("M134").Value = -4236.75
$ws.Range("H136").Value = 2906.4092
$ws.Range("I136").Value = 2763.8462
$ws.Range("J136").Value = 3112.3333
$ws.Range("K136").Value = 8291.5386
$ws.Range("L136").Value = 9336.999899999999
$ws.Range("M136").Value = -5741.5386
$ws.Range("N136").Value = -14436.9999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 1821.5454
$ws.Range("J29").Value = 2292.375
$ws.Range("L29").Value = 6877.125
$ws.Range("N29").Value = -7431.125
$ws.Range("H115").Value = 7762.25
$ws.Range("I115").Value = 8666.333000000001
$ws.Range("K115").Value = 25998.999
$ws.Range("M115").Value = -24823.999
$ws.Range("H138").Value = 7339326.5
$ws.Range("J138").Value = 8337908.5
$ws.Range("L138").Value = 25013725.5
$ws.Range("N138").Value = -25024005.5
$ws.Range("H140").Value = 3924.7
$ws.Range("I140").Value = 2869.2727
$ws.Range("K140").Value = 8607.8181
$ws.Range("M140").Value = -3427.8181

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 75998.5
$ws.Range("J62").Value = 75998.5
$ws.Range("L62").Value = 75998.5
$ws.Range("N62").Value = -77370.5
$ws.Range("H65").Value = 75998.5
$ws.Range("J65").Value = 75998.5
$ws.Range("L65").Value = 227995.5
$ws.Range("N65").Value = -234859.5
$ws.Range("H102").Value = 13719.768
$ws.Range("I102").Value = 1754.3429
$ws.Range("K102").Value = 1754.3429
$ws.Range("M102").Value = -132.3429000000001
$ws.Range("H107").Value = 746.5833
$ws.Range("I107").Value = 593.8461
$ws.Range("K107").Value = 593.8461
$ws.Range("M107").Value = 1326.1539
$ws.Range("H122").Value = 1849.1052
$ws.Range("I122").Value = 1592.6154
$ws.Range("K122").Value = 4777.8462
$ws.Range("M122").Value = -2327.8462
$ws.Range("H132").Value = 7094618
$ws.Range("I132").Value = 8132412
$ws.Range("J132").Value = 3024.6667
$ws.Range("K132").Value = 24397236
$ws.Range("L132").Value = 9074.000100000001
$ws.Range("M132").Value = -24394706
$ws.Range("N132").Value = -14134.0001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4090.3076
$ws.Range("I7").Value = 3462.121
$ws.Range("J7").Value = 5181.3687
$ws.Range("K7").Value = 3462.121
$ws.Range("L7").Value = 5181.3687
$ws.Range("M7").Value = -3350.121
$ws.Range("N7").Value = -5405.3687
$ws.Range("H40").Value = 8075.875
$ws.Range("I40").Value = 7855.7617
$ws.Range("K40").Value = 7855.7617
$ws.Range("M40").Value = -7719.7617
$ws.Range("H122").Value = 4177
$ws.Range("I122").Value = 3440.1482
$ws.Range("J122").Value = 5503.3335
$ws.Range("K122").Value = 10320.4446
$ws.Range("L122").Value = 16510.0005
$ws.Range("M122").Value = -7870.444600000001
$ws.Range("N122").Value = -21410.0005
$ws.Range("H126").Value = 4090.3076
$ws.Range("I126").Value = 3462.121
$ws.Range("J126").Value = 5181.3687
$ws.Range("K126").Value = 10386.363
$ws.Range("L126").Value = 15544.1061
$ws.Range("M126").Value = -7916.363000000001
$ws.Range("N126").Value = -20484.1061
$ws.Range("H132").Value = 3651.7344
$ws.Range("I132").Value = 3550
$ws.Range("K132").Value = 10650
$ws.Range("M132").Value = -8120
$ws.Range("H136").Value = 5818.4
$ws.Range("I136").Value = 4642.857
$ws.Range("J136").Value = 8561.333000000001
$ws.Range("K136").Value = 13928.571
$ws.Range("L136").Value = 25683.999
$ws.Range("M136").Value = -11378.571
$ws.Range("N136").Value = -30783.999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 32880.555
$ws.Range("J81").Value = 6548.706
$ws.Range("L81").Value = 13097.412
$ws.Range("N81").Value = -15219.412
$ws.Range("H84").Value = 32880.555
$ws.Range("J84").Value = 6548.706
$ws.Range("L84").Value = 65487.06
$ws.Range("N84").Value = -76095.06
$ws.Range("H107").Value = 1218.3572
$ws.Range("I107").Value = 1013.7
$ws.Range("J107").Value = 1730
$ws.Range("K107").Value = 3041.1
$ws.Range("L107").Value = 5190
$ws.Range("M107").Value = -1121.1
$ws.Range("N107").Value = -9030
$ws.Range("H122").Value = 4777.6
$ws.Range("I122").Value = 4777.6
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 14332.8
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -11882.8
$ws.Range("N122").ClearContents()
$ws.Range("H136").Value = 7760.8613
$ws.Range("I136").Value = 8602.620999999999
$ws.Range("K136").Value = 25807.863
$ws.Range("M136").Value = -23257.863

Write-Output "Applied 252 cell updates across 8 sheets"